$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-5 (Q0-Q3) with new values
$ws.Range("B2").Value = 0.03225293373050328
$ws.Range("C2").Value = 0.5338256249320711
$ws.Range("D2").Value = 0.4808498706638201
$ws.Range("E2").Value = 0.693433393098299
$ws.Range("F2").Value = 0.7003371009062773
$ws.Range("G2").Value = 46

$ws.Range("B3").Value = 0.08597907914849331
$ws.Range("C3").Value = 0.4451405877761389
$ws.Range("D3").Value = 0.430229207529593
$ws.Range("E3").Value = 0.655918598249503
$ws.Range("F3").Value = 0.6576068224065241
$ws.Range("G3").Value = 45

$ws.Range("B4").Value = 0.02586794070748817
$ws.Range("C4").Value = 0.4417395282174381
$ws.Range("D4").Value = 0.4101439193432779
$ws.Range("E4").Value = 0.6404247960090848
$ws.Range("F4").Value = 0.6473001143465348
$ws.Range("G4").Value = 44

$ws.Range("B5").Value = 0.06673040005454207
$ws.Range("C5").Value = 0.4545424064740749
$ws.Range("D5").Value = 0.4219317955494934
$ws.Range("E5").Value = 0.649562772601304
$ws.Range("F5").Value = 0.6537727601081559
$ws.Range("G5").Value = 43

# Add new rows 6-11 (Q4-Q9)
$ws.Range("A6").Value = "Q4"
$ws.Range("B6").Value = 0.07987673139996194
$ws.Range("C6").Value = 0.4232780668641379
$ws.Range("D6").Value = 0.4000004946887931
$ws.Range("E6").Value = 0.6324559231193847
$ws.Range("F6").Value = 0.6349966103947435
$ws.Range("G6").Value = 42

$ws.Range("A7").Value = "Q5"
$ws.Range("B7").Value = 0.05222626992688523
$ws.Range("C7").Value = 0.4466751155204058
$ws.Range("D7").Value = 0.4144800636829005
$ws.Range("E7").Value = 0.6438012610137546
$ws.Range("F7").Value = 0.649650900424786
$ws.Range("G7").Value = 41

$ws.Range("A8").Value = "Q6"
$ws.Range("B8").Value = 0.09078832835761171
$ws.Range("C8").Value = 0.4410619110573982
$ws.Range("D8").Value = 0.4058428363806655
$ws.Range("E8").Value = 0.6370579537064627
$ws.Range("F8").Value = 0.6385884400044998
$ws.Range("G8").Value = 40

$ws.Range("A9").Value = "Q7"
$ws.Range("B9").Value = 0.06725486063599991
$ws.Range("C9").Value = 0.4607570938983374
$ws.Range("D9").Value = 0.4414258711817555
$ws.Range("E9").Value = 0.6643988795759334
$ws.Range("F9").Value = 0.6696268312944889
$ws.Range("G9").Value = 39

$ws.Range("A10").Value = "Q8"
$ws.Range("B10").Value = 0.0892556709541302
$ws.Range("C10").Value = 0.4590231486236458
$ws.Range("D10").Value = 0.4345934021827165
$ws.Range("E10").Value = 0.6592369848413516
$ws.Range("F10").Value = 0.6619344999162985
$ws.Range("G10").Value = 38

$ws.Range("A11").Value = "Q9"
$ws.Range("B11").Value = 0.04838106079093066
$ws.Range("C11").Value = 0.4610557988666878
$ws.Range("D11").Value = 0.4449887869852402
$ws.Range("E11").Value = 0.6670747986434806
$ws.Range("F11").Value = 0.6744952478593288
$ws.Range("G11").Value = 37

# Apply the same style as A2:A5 (border + bold + alignment) to the new A6:A11 cells
$ws.Range("A2:A5").Copy()
$ws.Range("A6:A11").PasteSpecial(-4122)  # xlPasteFormats
